$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style = $ws.Range('D2').Style
$ws.Range('D2').Value = "'72.477.66"
$ws.Range('D2').Style = $style
$ws.Range('E2').Value = '  +4.33%  '

$style = $ws.Range('D3').Style
$ws.Range('D3').Value = "'4.049.90"
$ws.Range('D3').Style = $style
$ws.Range('E3').Value = '  +3.58%  '

$ws.Range('E4').Value = '  +0.06%  '

$style = $ws.Range('D5').Style
$ws.Range('D5').Value = "'519.40"
$ws.Range('D5').Style = $style
$ws.Range('E5').Value = '  -2.13%  '

$style = $ws.Range('D6').Style
$ws.Range('D6').Value = "'146.77"
$ws.Range('D6').Style = $style
$ws.Range('E6').Value = '  +1.26%  '

$ws.Range('E7').Value = '  +17.68%  '

$style = $ws.Range('D8').Style
$ws.Range('D8').Value = "'1.00"
$ws.Range('D8').Style = $style
$ws.Range('E8').Value = '  +0.12%  '

$style = $ws.Range('D9').Style
$ws.Range('D9').Value = "'0.758"
$ws.Range('D9').Style = $style
$ws.Range('E9').Value = '  +5.35%  '

$ws.Range('E10').Value = '  +1.10%  '

$ws.Range('E11').Value = '  -2.50%  '

$style = $ws.Range('D12').Style
$ws.Range('D12').Value = "'47.06"
$ws.Range('D12').Style = $style
$ws.Range('E12').Value = '  +11.51%  '

$style = $ws.Range('D13').Style
$ws.Range('D13').Value = "'10.89"
$ws.Range('D13').Style = $style
$ws.Range('E13').Value = '  +6.10%  '

$style = $ws.Range('D14').Style
$ws.Range('D14').Value = "'4.688.74"
$ws.Range('D14').Style = $style
$ws.Range('E14').Value = '  +3.50%  '

$style = $ws.Range('D15').Style
$ws.Range('D15').Value = "'4.037.81"
$ws.Range('D15').Style = $style
$ws.Range('E15').Value = '  +3.28%  '

$style = $ws.Range('D16').Style
$ws.Range('D16').Value = "'21.07"
$ws.Range('D16').Style = $style
$ws.Range('E16').Value = '  +6.43%  '

$style = $ws.Range('D17').Style
$ws.Range('D17').Value = "'14.13"
$ws.Range('D17').Style = $style
$ws.Range('E17').Value = '  +0.47%  '

$ws.Range('E18').Value = '  -1.71%  '

$ws.Range('E19').Value = '  -1.82%  '

$style = $ws.Range('D20').Style
$ws.Range('D20').Value = "'72.346.92"
$ws.Range('D20').Style = $style
$ws.Range('E20').Value = '  +4.23%  '

$style = $ws.Range('D21').Style
$ws.Range('D21').Value = "'442.55"
$ws.Range('D21').Style = $style
$ws.Range('E21').Value = '  +2.94%  '

$style = $ws.Range('D22').Style
$ws.Range('D22').Value = "'104.80"
$ws.Range('D22').Style = $style
$ws.Range('E22').Value = '  +18.34%  '

$style = $ws.Range('D23').Style
$ws.Range('D23').Value = "'3.59"
$ws.Range('D23').Style = $style
$ws.Range('E23').Value = '  +5.74%  '

$style = $ws.Range('D24').Style
$ws.Range('D24').Value = "'14.62"
$ws.Range('D24').Style = $style
$ws.Range('E24').Value = '  +2.83%  '

$style = $ws.Range('D25').Style
$ws.Range('D25').Value = "'3.99"
$ws.Range('D25').Style = $style
$ws.Range('E25').Value = '  -0.86%  '

$style = $ws.Range('D26').Style
$ws.Range('D26').Value = "'11.45"
$ws.Range('D26').Style = $style
$ws.Range('E26').Value = '  -0.75%  '

$style = $ws.Range('D27').Style
$ws.Range('D27').Value = "'11.05"
$ws.Range('D27').Style = $style
$ws.Range('E27').Value = '  +4.00%  '

$style = $ws.Range('D28').Style
$ws.Range('D28').Value = "'37.66"
$ws.Range('D28').Style = $style
$ws.Range('E28').Value = '  +3.35%  '

$style = $ws.Range('D29').Style
$ws.Range('D29').Value = "'5.83"
$ws.Range('D29').Style = $style
$ws.Range('E29').Value = '  +2.58%  '

$style = $ws.Range('D30').Style
$ws.Range('D30').Value = "'3.12"
$ws.Range('D30').Style = $style
$ws.Range('E30').Value = '  +10.15%  '

$style = $ws.Range('D31').Style
$ws.Range('D31').Value = "'13.60"
$ws.Range('D31').Style = $style
$ws.Range('E31').Value = '  +3.10%  '

$ws.Range('E32').Value = '  +2.68%  '

$style = $ws.Range('D33').Style
$ws.Range('D33').Value = "'677.69"
$ws.Range('D33').Style = $style
$ws.Range('E33').Value = '  -1.63%  '

$style = $ws.Range('D34').Style
$ws.Range('D34').Value = "'6.83"
$ws.Range('D34').Style = $style
$ws.Range('E34').Value = '  +13.95%  '

$style = $ws.Range('D35').Style
$ws.Range('D35').Value = "'67.40"
$ws.Range('D35').Style = $style
$ws.Range('E35').Value = '  -1.20%  '

$style = $ws.Range('D36').Style
$ws.Range('D36').Value = "'42.62"
$ws.Range('D36').Style = $style
$ws.Range('E36').Value = '  +6.33%  '

$ws.Range('E37').Value = '  -2.93%  '

$style = $ws.Range('D38').Style
$ws.Range('D38').Value = "'0.0₃0861"
$ws.Range('D38').Style = $style
$ws.Range('E38').Value = '  +1.61%  '

$style = $ws.Range('D39').Style
$ws.Range('D39').Value = "'3.59"
$ws.Range('D39').Style = $style
$ws.Range('E39').Value = '  +12.79%  '

$ws.Range('E40').Value = '  +0.64%  '

$ws.Range('E41').Value = '  +0.00%  '

$style = $ws.Range('D42').Style
$ws.Range('D42').Value = "'0.0495"
$ws.Range('D42').Style = $style
$ws.Range('E42').Value = '  +2.93%  '

$style = $ws.Range('D43').Style
$ws.Range('D43').Value = "'0.998"
$ws.Range('D43').Style = $style
$ws.Range('E43').Value = '  -0.25%  '

$style = $ws.Range('D44').Style
$ws.Range('D44').Value = "'3.23"
$ws.Range('D44').Style = $style
$ws.Range('E44').Value = '  +3.18%  '

$ws.Range('E45').Value = '  +12.34%  '

$ws.Range('E46').Value = '  -2.73%  '

$style = $ws.Range('D47').Style
$ws.Range('D47').Value = "'3.45"
$ws.Range('D47').Style = $style
$ws.Range('E47').Value = '  +2.62%  '

$ws.Range('E48').Value = '  +2.11%  '

$style = $ws.Range('D49').Style
$ws.Range('D49').Value = "'9.04"
$ws.Range('D49').Style = $style
$ws.Range('E49').Value = '  +6.43%  '

$ws.Range('E50').Value = '  +1.43%  '

$ws.Range('E51').Value = '  +0.91%  '

